$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the cached "last generated" date/time placeholders (the deck was
#    regenerated on 2020-08-04, previously cached as 2020-07-31) across every
#    master / layout that carries one.
# ---------------------------------------------------------------------------

# Handout master date placeholder (cs-CZ, "31.07.2020" style)
$p.HandoutMaster.Shapes.Item("Date Placeholder 2").TextFrame.TextRange.Text = "04.08.2020"

# Notes master date placeholder (en-GB, "31/07/2020" style)
$p.NotesMaster.Shapes.Item("Date Placeholder 2").TextFrame.TextRange.Text = "04/08/2020"

# Slide master date placeholder (en-CA, "2020-07-31" style)
$p.SlideMaster.Shapes.Item("Date Placeholder 3").TextFrame.TextRange.Text = "2020-08-04"

# "Blank" custom layout (slideLayout2) date placeholder (en-CA, "2020-07-31" style)
$blankLayout = $p.SlideMaster.CustomLayouts.Item(2)
$blankLayout.Shapes.Item("Date Placeholder 1").TextFrame.TextRange.Text = "2020-08-04"

# ---------------------------------------------------------------------------
# 2) Reposition shapes on slide 1 ("baseline" report) and slide 2
#    ("comparison" report) to match the regenerated report layout.
#    PowerPoint COM positions are expressed in points (1 pt = 914400/72 EMU).
# ---------------------------------------------------------------------------

$s1 = $p.Slides.Item(1)
$s2 = $p.Slides.Item(2)

# --- Slide 1 -----------------------------------------------------------
$s1.Shapes.Item("TextBox 4").Left = 572.4225196850393
$s1.Shapes.Item("TextBox 9").Left = 572.4225196850393
$s1.Shapes.Item("TextBox 13").Left = 572.4224409448819

# --- Slide 2 -----------------------------------------------------------
$tb13 = $s2.Shapes.Item("TextBox 13")
$tb13.Left = 451.6774015748031
$tb13.Top = 290.04653543307086

$rect2 = $s2.Shapes.Item("Rectangle 2")
$rect2.Left = 17.007874015748033
$rect2.Width = 799.9063779527559

$rect5 = $s2.Shapes.Item("Rectangle 5")
$rect5.Left = 444.96937007874016
$rect5.Top = 372.9493700787402

$rect6 = $s2.Shapes.Item("Rectangle 6")
$rect6.Left = 444.96937007874016
$rect6.Top = 451.2104724409449
